$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new row, carrying down the number formatting
# that is already applied to the row above (row 3) so the new B/C cells
# pick up the same cell style as the rest of the data.
$ws.Rows("3:3").Copy() | Out-Null
$ws.Rows("4:4").Insert(-4121, 0) | Out-Null  # xlShiftDown, xlFormatFromLeftOrAbove

# Fill in the new row's values
$ws.Range("A4").Value = "1FMCU9GD1HUA30879"
$ws.Range("B4").Value = 18000
$ws.Range("C4").Value = 25000

# Update the active selection to the new cell
$ws.Range("A4").Select()
